$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "End Km" column (column G, header + data) without shifting other columns
$ws.Range("G1:G16").ClearContents()

# Update the selected cell to match the new active selection
$ws.Range("F1").Select()
